# Reproduce the upstream edit: drop the leading "day of week" column
# (shifting everything one column to the left), refresh several of the
# "Responsible" / task names, and append two new rows of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove column A ("день недели") entirely - this shifts B:I -> A:H,
#    renumbers the merged header cells and keeps every other cell's
#    value/style attached to its (now shifted) address automatically.
$ws.Range("A1").EntireColumn.Delete()

# 2) Update the "Responsible" names / task text that changed in row 3 & 4.
$ws.Range("B3").Value = "Копытов П.Е.1"
$ws.Range("D3").Value = "Копытов П.Е.3"
$ws.Range("F3").Value = "Копытоа П.Е.5"

$ws.Range("B4").Value = "Копытов П.Е.2"
$ws.Range("C4").Value = "доработка парсера123"
$ws.Range("D4").Value = "цуац4"
$ws.Range("F4").Value = "Копытоа П.Е.6"

# 3) Row 5 (E5:F5) gets new content.
$ws.Range("E5").Value = "йуа"
$ws.Range("F5").Value = "12у"

# 4) Add a brand-new row 6 (E6:F6).
$ws.Range("E6").Value = "тест"
$ws.Range("F6").Value = "12ув"

# 5) Selection ends up on F6 in the saved file, and the frozen/scrolled
#    "topLeftCell" from the original view is no longer present.
$ws.Range("F6").Select()
